$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.210.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "'1.813.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.55%  "
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").Value = "'329.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.4444"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.09%  "
$ws.Range("D8").Value = "'0.3728"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("D9").Value = "'44.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'0.07696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.26%  "
$ws.Range("D11").Value = "'1.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "'22.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "'6.275"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "'7.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.60%  "
$ws.Range("D16").Value = "'1.818.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "'93.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +11.19%  "
$ws.Range("D18").Value = "'0.00001081"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("D19").Value = "'0.06520"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.51%  "
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'17.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.15%  "
$ws.Range("D22").Value = "'6.241"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.25%  "
$ws.Range("D23").Value = "'0.5338"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "'28.275.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").Value = "'11.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.43%  "
$ws.Range("D26").Value = "'2.064"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -14.32%  "
$ws.Range("D27").Value = "'20.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.30%  "
$ws.Range("D28").Value = "'154.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").Value = "'2.021.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.96%  "
$ws.Range("D30").Value = "'2.320"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'127.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").Value = "'1.198"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("D33").Value = "'5.856"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.52%  "
$ws.Range("D34").Value = "'0.09205"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").Value = "'3.677"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("E36").Value = "  +6.87%  "
$ws.Range("D37").Value = "'0.02347"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").Value = "'0.2170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").Value = "'5.176"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("D40").Value = "'0.06202"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").Value = "'0.6561"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").Value = "'1.199"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").Value = "'13.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("D46").Value = "'1.384"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "'0.6072"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("D48").Value = "'3.760"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "'126.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("D50").Value = "'2.033"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.98%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06982"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
